$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format, font, border, fill) from column G to column H
# for rows 4-37, matching the style each G cell uses in the same row.
$ws.Range("G4:G37").Copy() | Out-Null
$ws.Range("H4:H37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new 2022 data column (H).
$ws.Cells.Item(4, 8).Value2 = 2022
$ws.Cells.Item(5, 8).Value2 = 92.960099223795225
$ws.Cells.Item(7, 8).Value2 = 96.03949422949897
$ws.Cells.Item(8, 8).Value2 = 91.012153547624152
$ws.Cells.Item(10, 8).Value2 = 94.391087218067838
$ws.Cells.Item(11, 8).Value2 = 91.76755842559642
$ws.Cells.Item(13, 8).Value2 = 92.942689638142156
$ws.Cells.Item(14, 8).Value2 = 86.897877953385489
$ws.Cells.Item(15, 8).Value2 = 96.500794494289821
$ws.Cells.Item(16, 8).Value2 = 94.135975315309977
$ws.Cells.Item(17, 8).Value2 = 89.456106196597958
$ws.Cells.Item(18, 8).Value2 = 94.270923428904894
$ws.Cells.Item(19, 8).Value2 = 97.027480110114013
$ws.Cells.Item(20, 8).Value2 = 98.077227596867303
$ws.Cells.Item(21, 8).Value2 = 90.983384827072243
$ws.Cells.Item(23, 8).Value2 = 90.468970496790078
$ws.Cells.Item(24, 8).Value2 = 95.809965597614095
$ws.Cells.Item(25, 8).Value2 = 88.221110530662017
$ws.Cells.Item(27, 8).Value2 = 69.811292606515579
$ws.Cells.Item(28, 8).Value2 = 85.757158930558518
$ws.Cells.Item(29, 8).Value2 = 93.032103866435918
$ws.Cells.Item(30, 8).Value2 = 97.325262246493097
$ws.Cells.Item(31, 8).Value2 = 98.908492141713779
$ws.Cells.Item(33, 8).Value2 = 91.968006037496949
$ws.Cells.Item(34, 8).Value2 = 91.809335747904541
$ws.Cells.Item(35, 8).Value2 = 91.27524653351901
$ws.Cells.Item(36, 8).Value2 = 94.397288657466234
$ws.Cells.Item(37, 8).Value2 = 96.740699993405215

# Restore the active selection to I4, matching the saved view state.
$ws.Range("I4").Select() | Out-Null
